# Insert a new weekly price record as row 261, pushing the existing
# rows 261-265 down to 262-266.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 261 (shifts 261:265 down to 262:266).
$ws.Rows.Item(261).Insert()

# Populate the new row 261 with the new record. Columns A, B, C, E, F, G, R
# carry the same constant values as the surrounding rows in this block.
$ws.Range("A261").Value = 10
$ws.Range("B261").Value = "Vega Modelo de Temuco"
$ws.Range("C261").Value = "La Araucanía"
$ws.Range("D261").Value = 44448
$ws.Range("E261").Value = 9
$ws.Range("F261").Value = 100112024
$ws.Range("G261").Value = "Choclo"
$ws.Range("H261").Value = "Dulce o Americano"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 15
$ws.Range("K261").Value = 32000
$ws.Range("L261").Value = 32000
$ws.Range("M261").Value = 32000
$ws.Range("N261").Value = "$/malla 50 unidades"
$ws.Range("O261").Value = "Argentina"
$ws.Range("P261").Value = 640
$ws.Range("Q261").Value = 50
$ws.Range("R261").Value = "Hortaliza"
